# Updated cryptos list on Thu Mar  2 19:56:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and Column E (Volume(1h)) updates, keyed by row number.
# Only rows whose D value changed are listed in $prices; all rows with an
# E (percentage) change are listed in $volumes.

$prices = @{
  2  = "23.427.17"
  3  = "1.646.00"
  4  = "0.9996"
  5  = "0.9993"
  6  = "299.89"
  7  = "0.3800"
  9  = "0.3491"
  10 = "0.08070"
  12 = "0.9996"
  14 = "6.314"
  15 = "7.260"
  16 = "0.00001216"
  17 = "1.638.51"
  18 = "94.81"
  19 = "0.06959"
  20 = "6.628"
  22 = "0.9985"
  23 = "12.41"
  24 = "23.459.94"
  25 = "2.431"
  28 = "150.02"
  29 = "5.175"
  30 = "131.62"
  31 = "1.828.28"
  32 = "6.896"
  33 = "2.130"
  34 = "11.22"
  35 = "0.9902"
  36 = "0.02688"
  38 = "5.910"
  39 = "0.2422"
  40 = "0.06829"
  41 = "12.77"
  43 = "1.291"
  44 = "15.50"
  45 = "0.9980"
  46 = "0.6354"
  47 = "2.242"
  48 = "3.914"
  49 = "0.07678"
  50 = "127.04"
  51 = "1.225"
}

$volumes = @{
  2  = "  -0.46%  "
  3  = "  +0.02%  "
  4  = "  -0.65%  "
  5  = "  -0.65%  "
  6  = "  -1.10%  "
  7  = "  -0.69%  "
  8  = "  -1.63%  "
  9  = "  -2.94%  "
  10 = "  -1.20%  "
  11 = "  -0.29%  "
  12 = "  -0.65%  "
  13 = "  -1.14%  "
  14 = "  -1.70%  "
  15 = "  -1.89%  "
  16 = "  +0.01%  "
  17 = "  -0.84%  "
  18 = "  -2.58%  "
  19 = "  -1.04%  "
  20 = "  -1.91%  "
  21 = "  -0.40%  "
  22 = "  -0.61%  "
  23 = "  -1.85%  "
  24 = "  -0.31%  "
  25 = "  -1.88%  "
  26 = "  -1.62%  "
  27 = "  -0.94%  "
  28 = "  -2.19%  "
  29 = "  -1.10%  "
  30 = "  -1.80%  "
  31 = "  -0.49%  "
  32 = "  -2.08%  "
  33 = "  -5.86%  "
  34 = "  -6.79%  "
  35 = "  -6.00%  "
  36 = "  -3.79%  "
  37 = "  +0.16%  "
  38 = "  -2.04%  "
  39 = "  -2.82%  "
  40 = "  -1.97%  "
  41 = "  -1.38%  "
  42 = "  -1.82%  "
  43 = "  -3.14%  "
  44 = "  -2.82%  "
  45 = "  -0.56%  "
  46 = "  -1.77%  "
  47 = "  -2.11%  "
  48 = "  -1.26%  "
  49 = "  -2.21%  "
  50 = "  -0.44%  "
  51 = "  +2.41%  "
}

# These Price strings are digit-only (e.g. "0.9996"), so a plain .Value
# assignment would be auto-coerced to a number by Excel. The source data
# keeps them as literal text (inlineStr) in every row, so force the cell
# to Text format first, matching how the sheet was originally authored.
foreach ($row in $prices.Keys) {
  $cell = $ws.Cells.Item($row, 4)
  $cell.NumberFormat = "@"
  $cell.Value = $prices[$row]
}

foreach ($row in $volumes.Keys) {
  $ws.Cells.Item($row, 5).Value = $volumes[$row]
}

# Rows 38/39 swap places: Algorand <-> InternetComputer(DFINITY)
$ws.Cells.Item(38, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"

$ws.Cells.Item(39, 2).Value = "Algorand"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"

# Rows 49/50 swap places: Quant <-> Cronos
$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"

$ws.Cells.Item(50, 2).Value = "Quant"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
